$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AG2").Value = "MAB_and_AFACT"
$ws.Range("AJ2").Value = 7
$ws.Range("AG3").Value = "Dichotic_and_AFACT"
$ws.Range("AJ3").Value = 2
$ws.Range("AJ4").Value = 8
$ws.Range("AG5").Value = "MAB_phase"
$ws.Range("AJ5").Value = 8
$ws.Range("AG6").Value = "dichotic_phase"
$ws.Range("AG7").Value = "dichotic_phase"
$ws.Range("AJ7").Value = 2
$ws.Range("AG8").Value = "MAB_phase"
$ws.Range("AJ8").Value = 5
$ws.Range("AG9").Value = "MAB_phase"
$ws.Range("AJ9").Value = 1
$ws.Range("AG10").Value = "Dichotic_and_AFACT"
$ws.Range("AJ10").Value = 1
$ws.Range("AG11").Value = "dichotic_phase"
$ws.Range("AG12").Value = "Digit_before_and_AFACT"
$ws.Range("AJ12").Value = 7
$ws.Range("AG13").Value = "Digit_before_and_AFACT"
$ws.Range("AJ13").Value = 5
$ws.Range("AG14").Value = "dichotic_phase"
$ws.Range("AJ14").Value = 2
$ws.Range("AJ15").Value = 3
$ws.Range("AG16").Value = "MAB_and_Digit_after"
$ws.Range("AJ16").Value = 3
$ws.Range("AG17").Value = "Digit_before_and_AFACT"
$ws.Range("AJ17").Value = 7
$ws.Range("AG18").Value = "MAB_and_AFACT"
$ws.Range("AJ18").Value = 7
$ws.Range("AJ19").Value = 5
$ws.Range("AG20").Value = "Dichotic_and_AFACT"
$ws.Range("AJ20").Value = 4
$ws.Range("AJ21").Value = 6
$ws.Range("AG22").Value = "MAB_and_Digit_after"
$ws.Range("AJ22").Value = 2
$ws.Range("AJ23").Value = 2
$ws.Range("AG24").Value = "MAB_and_AFACT"
$ws.Range("AJ24").Value = 7
$ws.Range("AG25").Value = "Dichotic_and_AFACT"
$ws.Range("AJ25").Value = 7
$ws.Range("AG26").Value = "dichotic_phase"
$ws.Range("AJ26").Value = 2
$ws.Range("AG27").Value = "dichotic_phase"
$ws.Range("AJ27").Value = 6
$ws.Range("AJ28").Value = 1
$ws.Range("AG29").Value = "dichotic_phase"
$ws.Range("AJ29").Value = 1
$ws.Range("AJ30").Value = 2
$ws.Range("AG31").Value = "MAB_and_Digit_after"
$ws.Range("AJ31").Value = 5
$ws.Range("AG32").Value = "MAB_and_AFACT"
$ws.Range("AJ32").Value = 4
$ws.Range("AG33").Value = "Dichotic_and_AFACT"
$ws.Range("AJ33").Value = 3
$ws.Range("AG34").Value = "Digit_before_and_AFACT"
$ws.Range("AJ34").Value = 8
$ws.Range("AG35").Value = "Digit_before_and_AFACT"
$ws.Range("AJ35").Value = 3
$ws.Range("AJ36").Value = 3
$ws.Range("AG37").Value = "MAB_and_Digit_after"
$ws.Range("AJ37").Value = 3
$ws.Range("AJ38").Value = 8
$ws.Range("AJ39").Value = 3
$ws.Range("AG40").Value = "Dichotic_and_AFACT"
$ws.Range("AJ40").Value = 1
$ws.Range("AG41").Value = "MAB_phase"
$ws.Range("AG43").Value = "dichotic_phase"
$ws.Range("AJ43").Value = 6
$ws.Range("AG44").Value = "MAB_and_Digit_after"
$ws.Range("AJ44").Value = 6
$ws.Range("AG45").Value = "dichotic_phase"
$ws.Range("AJ45").Value = 5
$ws.Range("AG46").Value = "Dichotic_and_AFACT"
$ws.Range("AJ46").Value = 2
$ws.Range("AJ47").Value = 6
$ws.Range("AG48").Value = "Dichotic_and_AFACT"
$ws.Range("AJ48").Value = 3
$ws.Range("AJ49").Value = 1
$ws.Range("AG50").Value = "Dichotic_and_AFACT"
$ws.Range("AJ50").Value = 4
$ws.Range("AG51").Value = "dichotic_phase"
$ws.Range("AJ51").Value = 8
$ws.Range("AG52").Value = "MAB_phase"
$ws.Range("AJ52").Value = 4
$ws.Range("AG53").Value = "MAB_and_Digit_after"
$ws.Range("AJ53").Value = 1
$ws.Range("AG54").Value = "MAB_and_Digit_after"
$ws.Range("AJ54").Value = 5
$ws.Range("AG55").Value = "MAB_and_AFACT"
$ws.Range("AG56").Value = "Digit_before_and_AFACT"
$ws.Range("AJ56").Value = 6
$ws.Range("AJ57").Value = 3
$ws.Range("AG58").Value = "dichotic_phase"
$ws.Range("AJ58").Value = 7
$ws.Range("AG59").Value = "Digit_before_and_AFACT"
$ws.Range("AJ59").Value = 3
$ws.Range("AG60").Value = "MAB_and_AFACT"
$ws.Range("AJ60").Value = 5
$ws.Range("AG61").Value = "Digit_before_and_AFACT"
$ws.Range("AJ61").Value = 2
$ws.Range("AG62").Value = "MAB_phase"
$ws.Range("AJ62").Value = 8
$ws.Range("AG63").Value = "MAB_and_AFACT"
$ws.Range("AJ63").Value = 3
$ws.Range("AG64").Value = "MAB_phase"
$ws.Range("AJ64").Value = 1
$ws.Range("AG65").Value = "MAB_and_AFACT"
$ws.Range("AJ65").Value = 6
$ws.Range("AG66").Value = "dichotic_phase"
$ws.Range("AJ66").Value = 2
$ws.Range("AG67").Value = "Dichotic_and_AFACT"
$ws.Range("AG68").Value = "dichotic_phase"
$ws.Range("AJ69").Value = 7
$ws.Range("AJ70").Value = 5
$ws.Range("AG71").Value = "dichotic_phase"
$ws.Range("AJ71").Value = 3
$ws.Range("AG72").Value = "MAB_and_Digit_after"
$ws.Range("AJ72").Value = 4
$ws.Range("AG73").Value = "Dichotic_and_AFACT"
$ws.Range("AJ73").Value = 6
$ws.Range("AJ74").Value = 5
$ws.Range("AG75").Value = "Dichotic_and_AFACT"
$ws.Range("AJ75").Value = 6
$ws.Range("AG77").Value = "dichotic_phase"
$ws.Range("AJ77").Value = 2
$ws.Range("AG78").Value = "MAB_phase"
$ws.Range("AJ78").Value = 1
$ws.Range("AJ79").Value = 4
$ws.Range("AG80").Value = "Dichotic_and_AFACT"
$ws.Range("AJ80").Value = 7
$ws.Range("AG81").Value = "dichotic_phase"
$ws.Range("AJ81").Value = 8
